$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values
$ws.Range("G2").Value = 1
$ws.Range("J3").Value = 2
$ws.Range("M3").Value = 2
$ws.Range("J4").Value = 1
$ws.Range("M4").Value = 2

# Update the selection / active cell shown in the sheet view
$ws.Range("D3").Select()
